# [FIX] Se arreglaron casos para los modulos VDF y RC
# Adds 4 new user rows to the "Users" sheet (3rd worksheet):
#   F00481 / 081
#   F02214 / 221
#   F00197 / 007
#   F00042 / 042

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$ws.Activate()

# Row 62: F00481 / 081  (Sucursal kept as text, right aligned like existing rows)
$ws.Range("A62").Value = "F00481"
$ws.Range("C62").NumberFormat = "@"
$ws.Range("C62").HorizontalAlignment = -4152
$ws.Range("C62").Value = "081"

# Row 63: F02214 / 221 (numeric Sucursal, like rows 23/35/36/38/41/51/52/54)
$ws.Range("A63").Value = "F02214"
$ws.Range("C63").Value = 221

# Row 64: F00197 / 007 (text Sucursal)
$ws.Range("A64").Value = "F00197"
$ws.Range("C64").NumberFormat = "@"
$ws.Range("C64").HorizontalAlignment = -4152
$ws.Range("C64").Value = "007"

# Row 65: F00042 / 042 (text Sucursal)
$ws.Range("A65").Value = "F00042"
$ws.Range("C65").NumberFormat = "@"
$ws.Range("C65").HorizontalAlignment = -4152
$ws.Range("C65").Value = "042"

# Restore the view state left behind by the manual edit session
[void]$ws.Range("H68").Select()
$excel.ActiveWindow.ScrollRow = 49
